$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-14 Saturday" "2026-02-15 Sunday"

Replace-Text "295÷9=" "313÷4="
Replace-Text "354÷9=" "443÷4="
Replace-Text "448÷7=" "249÷2="
Replace-Text "765÷7=" "496÷5="
Replace-Text "998÷5=" "162÷3="
Replace-Text "811÷7=" "881÷7="
Replace-Text "950÷6=" "704÷3="
Replace-Text "825÷6=" "804÷4="
Replace-Text "921÷5=" "794÷7="
Replace-Text "925÷7=" "323÷7="
Replace-Text "753÷5=" "861÷6="
Replace-Text "739÷8=" "506÷9="
Replace-Text "337÷2=" "525÷6="
Replace-Text "654÷5=" "249÷8="
Replace-Text "236÷5=" "129÷7="
Replace-Text "132÷7=" "698÷2="
Replace-Text "178÷5=" "629÷6="
Replace-Text "473÷3=" "955÷5="
Replace-Text "400÷7=" "917÷3="
Replace-Text "775÷4=" "736÷2="
Replace-Text "578÷3=" "198÷8="
Replace-Text "323÷5=" "259÷3="
Replace-Text "327÷5=" "386÷7="
Replace-Text "533÷7=" "127÷7="
Replace-Text "356÷4=" "992÷6="
